# Sprint task board template: add a {teamname} placeholder after "Team: "
# and drop the old hard-coded "CHex" team name run.
#
# The original paragraph holds two separate runs:
#   1) "Team: "   (bold, normal size)
#   2) "CHex"     (bold, size 28 -- larger "team name" run)
# The edit folds the placeholder into run (1)'s text and removes run (2)
# entirely, so the paragraph ends up with a single "Team: {teamname}" run
# that keeps run (1)'s formatting.

$d = $word.ActiveDocument

# Remove the old literal team name run ("CHex") completely.
$nameRange = $d.Content
$nameRange.Find.Execute("CHex")
$nameRange.Text = ""

# Insert the {teamname} placeholder right after "Team: ", merging into
# that run since the formatting is identical.
$labelRange = $d.Content
$labelRange.Find.Execute("Team: ")
$labelRange.Collapse(0)
$labelRange.InsertAfter("{teamname}")
